$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.278.84'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.629.99'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.11'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.26'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.73%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +1.07%  '
$ws.Range("E9").Value = '  +2.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.77'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.29%  '
$ws.Range("E11").Value = '  +6.94%  '
$ws.Range("E12").Value = '  -0.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.69'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.104.64'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.129.47'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.98%  '
$ws.Range("E16").Value = '  +4.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.627.85'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.17'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +8.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.67'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '349.82'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.02'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.74'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.72'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("E25").Value = '  +15.54%  '
$ws.Range("E26").Value = '  +5.65%  '
$ws.Range("E27").Value = '  +7.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.166'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.11'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '547.19'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.15%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("E32").Value = '  +2.59%  '
$ws.Range("E33").Value = '  +7.58%  '
$ws.Range("E34").Value = '  +0.90%  '
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '167.26'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("E37").Value = '  +8.21%  '
$ws.Range("E38").Value = '  +2.48%  '
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.52'
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '173.10'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.99'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.94'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +5.63%  '
$ws.Range("E45").Value = '  +4.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.64'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.31%  '
$ws.Range("E47").Value = '  +0.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.02'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +15.51%  '
$ws.Range("E49").Value = '  +2.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0967'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.31'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +4.36%  '
